$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param([int]$row1, [int]$row2)

    # Columns B through AD hold the match record (column A is just the
    # positional row id and must stay where it is).
    $rng1 = $ws.Range("B$row1" + ":AD$row1")
    $rng2 = $ws.Range("B$row2" + ":AD$row2")

    # NOTE: read via .Value2 -- the .Value getter does not reliably
    # return data when assigned into a variable in this host.
    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value = $vals2
    $rng2.Value = $vals1
}

# Rows 172 and 174 had their records swapped (positional args only --
# named parameter binding is unreliable for functions in this host).
Swap-Rows 172 174

# Rows 176 and 177 had their records swapped.
Swap-Rows 176 177
